# realize PlayerTank and PlayerController gameplay logic
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "输入：按键映射" (input: key mapping) is now done -> mark Good (was Neutral)
$ws.Range("A39").Style = "Good"

# Fix typo: 效果 (effect) -> 特效 (special effect) in the combined task summary cell
$ws.Range("A59").Value = "特效，碰撞，物理，声音，输入，材质，纹理"

# Add "completed this week" notes for the tank control / player controller rows
$ws.Range("B33").Value = "本周完成控制，瞄准"
$ws.Range("B35").Value = "本周完成"

# Break down the combined "菜单，本地化，加载，关卡" task into individual tracked items
$ws.Range("B58").Value = "菜单"
$ws.Range("B58").Style = "Good"
$ws.Range("C58").Value = "关卡"
$ws.Range("C58").Style = "Neutral"
$ws.Range("D58").Value = "本地化"
$ws.Range("D58").Style = "Bad"
$ws.Range("E58").Value = "加载"
$ws.Range("E58").Style = "Bad"

# Break down the combined "特效，碰撞，物理，声音，输入，材质，纹理" task into individual tracked items
$ws.Range("B59").Value = "特效"
$ws.Range("B59").Style = "Bad"
$ws.Range("C59").Value = "碰撞物理"
$ws.Range("C59").Style = "Good"
$ws.Range("D59").Value = "声音"
$ws.Range("D59").Style = "Good"
$ws.Range("E59").Value = "输入"
$ws.Range("E59").Style = "Good"
$ws.Range("F59").Value = "材质"
$ws.Range("F59").Style = "Good"
$ws.Range("G59").Value = "纹理"
$ws.Range("G59").Style = "Good"

# Move the active selection to match the author's final cursor position
$ws.Range("B35").Select()
